$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear row 18 data cells (B18:K18) - keep only the "sd" label in A18
$ws.Range("B18:K18").Clear()

# Update cell values with resubmission-recomputed figures (tiny floating point diffs)
$ws.Range("C2").Value = 0.65632594318798965
$ws.Range("E2").Value = 0.3305159220204571
$ws.Range("G2").Value = 0.25016840817816172
$ws.Range("H3").Value = 0.22784349420590899
$ws.Range("B4").Value = -0.30924910681999701
$ws.Range("C4").Value = 0.18766164446218619
$ws.Range("D4").Value = -0.53570588841842315
$ws.Range("E4").Value = -0.28069368237635067
$ws.Range("F4").Value = -0.31102196419501732
$ws.Range("G4").Value = -0.016654211424877399
$ws.Range("D5").Value = 0.13945749866997051
$ws.Range("E5").Value = 0.10613850178552971
$ws.Range("F5").Value = 0.11485691821640311
$ws.Range("G5").Value = 0.10698312800860239
$ws.Range("J5").Value = 0.099923330369449698
$ws.Range("K5").Value = 0.099297407939974103
$ws.Range("B8").Value = 0.33185103393838378
$ws.Range("G8").Value = 0.065174633976522206
$ws.Range("B10").Value = 0.63648765648375583
$ws.Range("D10").Value = 0.67669918391932216
$ws.Range("E10").Value = 0.061948692914857599
$ws.Range("H10").Value = -0.28599178792943958
$ws.Range("B11").Value = 0.092827431765593202
$ws.Range("C11").Value = 0.23206510832558219
$ws.Range("D11").Value = 0.0989524393338337
$ws.Range("E11").Value = 0.099048268081384394
$ws.Range("G11").Value = 0.19289050033005661
$ws.Range("J11").Value = 0.1979013036259662
$ws.Range("K11").Value = 0.2003731504081066
$ws.Range("E12").Value = 0.88682104894119662
$ws.Range("K12").Value = 0.54038233577596106
$ws.Range("F13").Value = 0.0630686022056023
$ws.Range("B14").Value = -0.64456029503611922
$ws.Range("D14").Value = -0.5196136538793833
$ws.Range("E14").Value = 0.1558603406596912
$ws.Range("F14").Value = 0.23206396498343521
$ws.Range("K14").Value = -0.081494607352841295
$ws.Range("B15").Value = 0.26397955583748139
$ws.Range("D15").Value = 0.25540428383690078
$ws.Range("E15").Value = 0.22208801137037459
$ws.Range("F15").Value = 0.2286130124949492
$ws.Range("G15").Value = 0.3070868420515907
$ws.Range("H15").Value = 0.97865895263828517
$ws.Range("I15").Value = 0.30089445442266671
$ws.Range("J15").Value = 0.2917033548022272
$ws.Range("K15").Value = 0.2933996004248593
$ws.Range("B19").Value = 0.083532623054420801
$ws.Range("C19").Value = 0.14663053038959711
$ws.Range("D19").Value = 0.060736480781608197
$ws.Range("E19").Value = 0.26062184677659389
$ws.Range("F19").Value = 0.24285765485841279
$ws.Range("G19").Value = 0.35306880904707377
$ws.Range("H19").Value = 0.47404065179798299
$ws.Range("I19").Value = 0.33969377760306441
$ws.Range("J19").Value = 0.37502216430269653
$ws.Range("K19").Value = 0.36249087069569103
